$d = $word.ActiveDocument

# Phase 1: create 4 placeholder paragraphs after the existing trailing empty
# paragraph (kept intact as the leading "`r" preserves it untouched).
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertAfter("`rZZPLACEHOLDER1ZZ`rZZPLACEHOLDER2ZZ`rZZPLACEHOLDER3ZZ`rZZPLACEHOLDER4ZZ")

# Phase 2.1: fill placeholder 1
$r0 = $d.Content
$r0.Find.Execute("ZZPLACEHOLDER1ZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r0.Text = ""
$xml0 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:r>
        <w:br w:type="page"/>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r0.InsertXML($xml0)

# Phase 2.2: fill placeholder 2
$r1 = $d.Content
$r1.Find.Execute("ZZPLACEHOLDER2ZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.Text = ""
$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Précisions sur le travail effectué</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r1.InsertXML($xml1)

# Phase 2.3: fill placeholder 3
$r2 = $d.Content
$r2.Find.Execute("ZZPLACEHOLDER3ZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Text = ""
$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r2.InsertXML($xml2)

# Phase 2.4: fill placeholder 4
$r3 = $d.Content
$r3.Find.Execute("ZZPLACEHOLDER4ZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Text = ""
$xml3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:r>
        <w:t xml:space="preserve">On peut voir que l’un de nous n’a pas fait beaucoup de </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>commits</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> sur le projet. Ceci est lié au fait qu’il a travaillé avec Arthur sur la fonctionnalité de la recherche des </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>pokémon</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">. Il a de plus, aidé Robin sur le rendu visuel de certaines pages (ex : détail du </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>pokemon</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> sélectionné).</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r3.InsertXML($xml3)

Write-Host "done"